# Fix bug caused by mutable default parameters so the "glass" sheet's
# rest-sector specific consumption values (Electricity / Heat) are
# calculated per-country instead of all sharing the same default values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("glass")

# Row => (Electricity [B], Heat [C])
$updates = @{
    2  = @(1.58, 7.48)   # Belgium
    3  = @(1.38, 8.34)   # Bulgaria
    5  = @(2.45, 4.62)   # Denmark
    7  = @(2.97, 8.94)   # Ireland
    8  = @(1.34, 5.71)   # Greece
    9  = @(1.29, 6.28)   # Spain
    10 = @(1.44, 6.75)   # France
    11 = @(2.15, 7.83)   # Croatia
    12 = @(1.13, 5.99)   # Italy
    13 = @(3.38, 10.32)  # Latvia
    15 = @(1.9,  7.02)   # Hungary
    18 = @(1.69, 6.22)   # Poland
    19 = @(1.1,  5.74)   # Portugal
    20 = @(1.18, 7.16)   # Romania
    21 = @(3.36, 4.31)   # Slovenia
    22 = @(2.36, 7.02)   # Slovakia
    23 = @(2.25, 6.4)    # Finland
    24 = @(1.91, 5.94)   # Sweden
    25 = @(1.28, 6.01)   # United Kingdom
    31 = @(0.99, 5.41)   # Serbia
    34 = @(3.38, 4.12)   # Lithuania
    35 = @(0.99, 5.41)   # Estonia
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
